# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Leviathan_Profits workbook (per-sheet Leve profit calculations).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2880.923  # H86: 2950.0908 -> 2880.923
$ws.Cells.Item(86, 9).Value = 2905.7778  # I86: 3049.5715 -> 2905.7778
$ws.Cells.Item(86, 10).Value = 2825  # J86: 2776 -> 2825
$ws.Cells.Item(86, 11).Value = 2905.7778  # K86: 3049.5715 -> 2905.7778
$ws.Cells.Item(86, 12).Value = 2825  # L86: 2776 -> 2825
$ws.Cells.Item(86, 13).Value = -1782.7778  # M86: -1926.5715 -> -1782.7778
$ws.Cells.Item(86, 14).Value = -5071  # N86: -5022 -> -5071
$ws.Cells.Item(89, 8).Value = 2880.923  # H89: 2950.0908 -> 2880.923
$ws.Cells.Item(89, 9).Value = 2905.7778  # I89: 3049.5715 -> 2905.7778
$ws.Cells.Item(89, 10).Value = 2825  # J89: 2776 -> 2825
$ws.Cells.Item(89, 11).Value = 14528.889  # K89: 15247.8575 -> 14528.889
$ws.Cells.Item(89, 12).Value = 14125  # L89: 13880 -> 14125
$ws.Cells.Item(89, 13).Value = -8912.888999999999  # M89: -9631.8575 -> -8912.888999999999
$ws.Cells.Item(89, 14).Value = -25357  # N89: -25112 -> -25357
$ws.Cells.Item(106, 8).Value = 14952.111  # H106: 15063.777 -> 14952.111
$ws.Cells.Item(106, 10).Value = 20997.666  # J106: 21165.166 -> 20997.666
$ws.Cells.Item(106, 12).Value = 20997.666  # L106: 21165.166 -> 20997.666
$ws.Cells.Item(106, 14).Value = -22259.666  # N106: -22427.166 -> -22259.666
$ws.Cells.Item(116, 8).Value = 5000  # H116: 3999 -> 5000
$ws.Cells.Item(116, 9).Value = 0  # I116: 2998 -> 0
$ws.Cells.Item(116, 11).Value = 0  # K116: 2998 -> 0
$ws.Cells.Item(116, 13).ClearContents()  # M116: 444 -> (removed)
$ws.Cells.Item(127, 8).Value = 910.1539  # H127: 1018.1539 -> 910.1539
$ws.Cells.Item(127, 9).Value = 893.8182  # I127: 1021.4545 -> 893.8182
$ws.Cells.Item(127, 11).Value = 2681.4546  # K127: 3064.3635 -> 2681.4546
$ws.Cells.Item(127, 13).Value = 2278.5454  # M127: 1895.6365 -> 2278.5454
$ws.Cells.Item(137, 8).Value = 5155.2  # H137: 4629.3335 -> 5155.2
$ws.Cells.Item(137, 10).Value = 7333.3335  # J137: 6000 -> 7333.3335
$ws.Cells.Item(137, 12).Value = 22000.0005  # L137: 18000 -> 22000.0005
$ws.Cells.Item(137, 14).Value = -27100.0005  # N137: -23100 -> -27100.0005
$ws.Cells.Item(138, 8).Value = 3629.6287  # H138: 3620.3333 -> 3629.6287
$ws.Cells.Item(138, 10).Value = 3768.6897  # J138: 3752.9 -> 3768.6897
$ws.Cells.Item(138, 12).Value = 11306.0691  # L138: 11258.7 -> 11306.0691
$ws.Cells.Item(138, 14).Value = -21586.0691  # N138: -21538.7 -> -21586.0691

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1623.9445  # H2: 1565.5264 -> 1623.9445
$ws.Cells.Item(2, 9).Value = 1535.1428  # I2: 1403.125 -> 1535.1428
$ws.Cells.Item(2, 10).Value = 1934.75  # J2: 2431.6667 -> 1934.75
$ws.Cells.Item(2, 11).Value = 1535.1428  # K2: 1403.125 -> 1535.1428
$ws.Cells.Item(2, 12).Value = 1934.75  # L2: 2431.6667 -> 1934.75
$ws.Cells.Item(2, 13).Value = -1422.1428  # M2: -1290.125 -> -1422.1428
$ws.Cells.Item(2, 14).Value = -2160.75  # N2: -2657.6667 -> -2160.75
$ws.Cells.Item(4, 8).Value = 126072.625  # H4: 126156.75 -> 126072.625
$ws.Cells.Item(4, 9).Value = 167589.25  # I4: 182779.27 -> 167589.25
$ws.Cells.Item(4, 10).Value = 1522.75  # J4: 1587.2 -> 1522.75
$ws.Cells.Item(4, 11).Value = 167589.25  # K4: 182779.27 -> 167589.25
$ws.Cells.Item(4, 12).Value = 1522.75  # L4: 1587.2 -> 1522.75
$ws.Cells.Item(4, 13).Value = -167473.25  # M4: -182663.27 -> -167473.25
$ws.Cells.Item(4, 14).Value = -1754.75  # N4: -1819.2 -> -1754.75
$ws.Cells.Item(32, 8).Value = 5801.844  # H32: 5731.4614 -> 5801.844
$ws.Cells.Item(32, 9).Value = 4121  # I32: 4062.4 -> 4121
$ws.Cells.Item(32, 11).Value = 4121  # K32: 4062.4 -> 4121
$ws.Cells.Item(32, 13).Value = -3834  # M32: -3775.4 -> -3834
$ws.Cells.Item(74, 8).Value = 1724.375  # H74: 1777.6086 -> 1724.375
$ws.Cells.Item(74, 9).Value = 924.6875  # I74: 953 -> 924.6875
$ws.Cells.Item(74, 11).Value = 924.6875  # K74: 953 -> 924.6875
$ws.Cells.Item(74, 13).Value = -50.6875  # M74: -79 -> -50.6875
$ws.Cells.Item(77, 8).Value = 1724.375  # H77: 1777.6086 -> 1724.375
$ws.Cells.Item(77, 9).Value = 924.6875  # I77: 953 -> 924.6875
$ws.Cells.Item(77, 11).Value = 4623.4375  # K77: 4765 -> 4623.4375
$ws.Cells.Item(77, 13).Value = -255.4375  # M77: -397 -> -255.4375
$ws.Cells.Item(116, 8).Value = 1623.9445  # H116: 1565.5264 -> 1623.9445
$ws.Cells.Item(116, 9).Value = 1535.1428  # I116: 1403.125 -> 1535.1428
$ws.Cells.Item(116, 10).Value = 1934.75  # J116: 2431.6667 -> 1934.75
$ws.Cells.Item(116, 11).Value = 1535.1428  # K116: 1403.125 -> 1535.1428
$ws.Cells.Item(116, 12).Value = 1934.75  # L116: 2431.6667 -> 1934.75
$ws.Cells.Item(116, 13).Value = 758.8571999999999  # M116: 890.875 -> 758.8571999999999
$ws.Cells.Item(116, 14).Value = -6522.75  # N116: -7019.6667 -> -6522.75
$ws.Cells.Item(122, 8).Value = 2392.4285  # H122: 2561.5417 -> 2392.4285
$ws.Cells.Item(122, 9).Value = 1701.9474  # I122: 1788.4 -> 1701.9474
$ws.Cells.Item(122, 11).Value = 5105.8422  # K122: 5365.200000000001 -> 5105.8422
$ws.Cells.Item(122, 13).Value = -2655.8422  # M122: -2915.200000000001 -> -2655.8422

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1623.9445  # H3: 1565.5264 -> 1623.9445
$ws.Cells.Item(3, 9).Value = 1535.1428  # I3: 1403.125 -> 1535.1428
$ws.Cells.Item(3, 10).Value = 1934.75  # J3: 2431.6667 -> 1934.75
$ws.Cells.Item(3, 11).Value = 1535.1428  # K3: 1403.125 -> 1535.1428
$ws.Cells.Item(3, 12).Value = 1934.75  # L3: 2431.6667 -> 1934.75
$ws.Cells.Item(3, 13).Value = -1421.1428  # M3: -1289.125 -> -1421.1428
$ws.Cells.Item(3, 14).Value = -2162.75  # N3: -2659.6667 -> -2162.75
$ws.Cells.Item(20, 10).Value = 4968.6665  # J20: 4969 -> 4968.6665
$ws.Cells.Item(20, 12).Value = 4968.6665  # L20: 4969 -> 4968.6665
$ws.Cells.Item(20, 14).Value = -5462.6665  # N20: -5463 -> -5462.6665
$ws.Cells.Item(86, 8).Value = 3214.5  # H86: 2877.2 -> 3214.5
$ws.Cells.Item(86, 9).Value = 2476.111  # I86: 2323.5454 -> 2476.111
$ws.Cells.Item(86, 10).Value = 5429.6665  # J86: 4399.75 -> 5429.6665
$ws.Cells.Item(86, 11).Value = 2476.111  # K86: 2323.5454 -> 2476.111
$ws.Cells.Item(86, 12).Value = 5429.6665  # L86: 4399.75 -> 5429.6665
$ws.Cells.Item(86, 13).Value = -1353.111  # M86: -1200.5454 -> -1353.111
$ws.Cells.Item(86, 14).Value = -7675.6665  # N86: -6645.75 -> -7675.6665
$ws.Cells.Item(89, 8).Value = 3214.5  # H89: 2877.2 -> 3214.5
$ws.Cells.Item(89, 9).Value = 2476.111  # I89: 2323.5454 -> 2476.111
$ws.Cells.Item(89, 10).Value = 5429.6665  # J89: 4399.75 -> 5429.6665
$ws.Cells.Item(89, 11).Value = 12380.555  # K89: 11617.727 -> 12380.555
$ws.Cells.Item(89, 12).Value = 27148.3325  # L89: 21998.75 -> 27148.3325
$ws.Cells.Item(89, 13).Value = -6764.555  # M89: -6001.726999999999 -> -6764.555
$ws.Cells.Item(89, 14).Value = -38380.3325  # N89: -33230.75 -> -38380.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 40689.81  # H31: 39253.52 -> 40689.81
$ws.Cells.Item(31, 9).Value = 54143.21  # I31: 51506.55 -> 54143.21
$ws.Cells.Item(31, 10).Value = 4173.4287  # J31: 4244.857 -> 4173.4287
$ws.Cells.Item(31, 11).Value = 54143.21  # K31: 51506.55 -> 54143.21
$ws.Cells.Item(31, 12).Value = 4173.4287  # L31: 4244.857 -> 4173.4287
$ws.Cells.Item(31, 13).Value = -53848.21  # M31: -51211.55 -> -53848.21
$ws.Cells.Item(31, 14).Value = -4763.4287  # N31: -4834.857 -> -4763.4287
$ws.Cells.Item(34, 8).Value = 40689.81  # H34: 39253.52 -> 40689.81
$ws.Cells.Item(34, 9).Value = 54143.21  # I34: 51506.55 -> 54143.21
$ws.Cells.Item(34, 10).Value = 4173.4287  # J34: 4244.857 -> 4173.4287
$ws.Cells.Item(34, 11).Value = 54143.21  # K34: 51506.55 -> 54143.21
$ws.Cells.Item(34, 12).Value = 4173.4287  # L34: 4244.857 -> 4173.4287
$ws.Cells.Item(34, 13).Value = -53941.21  # M34: -51304.55 -> -53941.21
$ws.Cells.Item(34, 14).Value = -4577.4287  # N34: -4648.857 -> -4577.4287
$ws.Cells.Item(51, 8).Value = 9750  # H51: 9000 -> 9750
$ws.Cells.Item(61, 8).Value = 9750  # H61: 9000 -> 9750
$ws.Cells.Item(68, 8).Value = 19424.285  # H68: 20746.25 -> 19424.285
$ws.Cells.Item(68, 9).Value = 10000  # I68: 0 -> 10000
$ws.Cells.Item(68, 10).Value = 20995  # J68: 20746.25 -> 20995
$ws.Cells.Item(68, 11).Value = 10000  # K68: 0 -> 10000
$ws.Cells.Item(68, 12).Value = 20995  # L68: 20746.25 -> 20995
$ws.Cells.Item(68, 13).Value = -9251  # M68: None -> -9251
$ws.Cells.Item(68, 14).Value = -22493  # N68: -22244.25 -> -22493
$ws.Cells.Item(69, 8).Value = 10250  # H69: 12000 -> 10250
$ws.Cells.Item(69, 9).Value = 10250  # I69: 12000 -> 10250
$ws.Cells.Item(69, 11).Value = 10250  # K69: 12000 -> 10250
$ws.Cells.Item(69, 13).Value = -9501  # M69: -11251 -> -9501
$ws.Cells.Item(71, 8).Value = 19424.285  # H71: 20746.25 -> 19424.285
$ws.Cells.Item(71, 9).Value = 10000  # I71: 0 -> 10000
$ws.Cells.Item(71, 10).Value = 20995  # J71: 20746.25 -> 20995
$ws.Cells.Item(71, 11).Value = 30000  # K71: 0 -> 30000
$ws.Cells.Item(71, 12).Value = 62985  # L71: 62238.75 -> 62985
$ws.Cells.Item(71, 13).Value = -26256  # M71: None -> -26256
$ws.Cells.Item(71, 14).Value = -70473  # N71: -69726.75 -> -70473
$ws.Cells.Item(72, 8).Value = 10250  # H72: 12000 -> 10250
$ws.Cells.Item(72, 9).Value = 10250  # I72: 12000 -> 10250
$ws.Cells.Item(72, 11).Value = 30750  # K72: 36000 -> 30750
$ws.Cells.Item(72, 13).Value = -27006  # M72: -32256 -> -27006

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 256.23077  # H23: 293.92307 -> 256.23077
$ws.Cells.Item(23, 9).Value = 285.6  # I23: 280.54544 -> 285.6
$ws.Cells.Item(23, 10).Value = 158.33333  # J23: 367.5 -> 158.33333
$ws.Cells.Item(23, 11).Value = 856.8000000000001  # K23: 841.63632 -> 856.8000000000001
$ws.Cells.Item(23, 12).Value = 474.99999  # L23: 1102.5 -> 474.99999
$ws.Cells.Item(23, 13).Value = -621.8000000000001  # M23: -606.63632 -> -621.8000000000001
$ws.Cells.Item(23, 14).Value = -944.99999  # N23: -1572.5 -> -944.99999
$ws.Cells.Item(68, 8).Value = 1285.3529  # H68: 1288.2941 -> 1285.3529
$ws.Cells.Item(68, 9).Value = 1064.4286  # I68: 1071.5714 -> 1064.4286
$ws.Cells.Item(68, 11).Value = 3193.2858  # K68: 3214.7142 -> 3193.2858
$ws.Cells.Item(68, 13).Value = -2382.2858  # M68: -2403.7142 -> -2382.2858
$ws.Cells.Item(71, 8).Value = 1285.3529  # H71: 1288.2941 -> 1285.3529
$ws.Cells.Item(71, 9).Value = 1064.4286  # I71: 1071.5714 -> 1064.4286
$ws.Cells.Item(71, 11).Value = 9579.857399999999  # K71: 9644.142600000001 -> 9579.857399999999
$ws.Cells.Item(71, 13).Value = -5523.857399999999  # M71: -5588.142600000001 -> -5523.857399999999
$ws.Cells.Item(80, 8).Value = 4916.6665  # H80: 5000 -> 4916.6665
$ws.Cells.Item(80, 10).Value = 4909.091  # J80: 5000 -> 4909.091
$ws.Cells.Item(80, 12).Value = 14727.273  # L80: 15000 -> 14727.273
$ws.Cells.Item(80, 14).Value = -16599.273  # N80: -16872 -> -16599.273
$ws.Cells.Item(83, 8).Value = 4916.6665  # H83: 5000 -> 4916.6665
$ws.Cells.Item(83, 10).Value = 4909.091  # J83: 5000 -> 4909.091
$ws.Cells.Item(83, 12).Value = 44181.819  # L83: 45000 -> 44181.819
$ws.Cells.Item(83, 14).Value = -53541.819  # N83: -54360 -> -53541.819
$ws.Cells.Item(139, 8).Value = 1881.2  # H139: 54234.95 -> 1881.2
$ws.Cells.Item(139, 9).Value = 1881.2  # I139: 54234.95 -> 1881.2
$ws.Cells.Item(139, 11).Value = 5643.6  # K139: 162704.85 -> 5643.6
$ws.Cells.Item(139, 13).Value = -503.6000000000004  # M139: -157564.85 -> -503.6000000000004

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 42682.363  # H7: 44357.76 -> 42682.363
$ws.Cells.Item(7, 9).Value = 49667.168  # I7: 58100.8 -> 49667.168
$ws.Cells.Item(7, 11).Value = 49667.168  # K7: 58100.8 -> 49667.168
$ws.Cells.Item(7, 13).Value = -49555.168  # M7: -57988.8 -> -49555.168
$ws.Cells.Item(18, 8).Value = 10000  # H18: 0 -> 10000
$ws.Cells.Item(18, 10).Value = 10000  # J18: 0 -> 10000
$ws.Cells.Item(18, 12).Value = 10000  # L18: 0 -> 10000
$ws.Cells.Item(18, 14).Value = -10344  # N18: None -> -10344
$ws.Cells.Item(40, 8).Value = 14379.333  # H40: 15424.643 -> 14379.333
$ws.Cells.Item(40, 9).Value = 7208.8096  # I40: 7994.4736 -> 7208.8096
$ws.Cells.Item(40, 11).Value = 7208.8096  # K40: 7994.4736 -> 7208.8096
$ws.Cells.Item(40, 13).Value = -7072.8096  # M40: -7858.4736 -> -7072.8096
$ws.Cells.Item(55, 8).Value = 1655.2222  # H55: 1786 -> 1655.2222
$ws.Cells.Item(55, 9).Value = 1537.125  # I55: 1786 -> 1537.125
$ws.Cells.Item(55, 10).Value = 2600  # J55: 0 -> 2600
$ws.Cells.Item(55, 11).Value = 1537.125  # K55: 1786 -> 1537.125
$ws.Cells.Item(55, 12).Value = 2600  # L55: 0 -> 2600
$ws.Cells.Item(55, 13).Value = -1364.125  # M55: -1613 -> -1364.125
$ws.Cells.Item(55, 14).Value = -2946  # N55: None -> -2946
$ws.Cells.Item(126, 8).Value = 42682.363  # H126: 44357.76 -> 42682.363
$ws.Cells.Item(126, 9).Value = 49667.168  # I126: 58100.8 -> 49667.168
$ws.Cells.Item(126, 11).Value = 149001.504  # K126: 174302.4 -> 149001.504
$ws.Cells.Item(126, 13).Value = -146531.504  # M126: -171832.4 -> -146531.504
$ws.Cells.Item(132, 8).Value = 4489.2905  # H132: 4488.933 -> 4489.2905
$ws.Cells.Item(132, 9).Value = 3819  # I132: 3791.76 -> 3819
$ws.Cells.Item(132, 11).Value = 11457  # K132: 11375.28 -> 11457
$ws.Cells.Item(132, 13).Value = -8927  # M132: -8845.280000000001 -> -8927

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 20000  # H29: 15000 -> 20000
$ws.Cells.Item(29, 9).Value = 0  # I29: 10000 -> 0
$ws.Cells.Item(29, 11).Value = 0  # K29: 10000 -> 0
$ws.Cells.Item(29, 13).ClearContents()  # M29: -9710 -> (removed)
$ws.Cells.Item(46, 8).Value = 91999.5  # H46: 77116.836 -> 91999.5
$ws.Cells.Item(46, 9).Value = 0  # I46: 70000 -> 0
$ws.Cells.Item(46, 10).Value = 91999.5  # J46: 78540.2 -> 91999.5
$ws.Cells.Item(46, 11).Value = 0  # K46: 70000 -> 0
$ws.Cells.Item(46, 12).Value = 91999.5  # L46: 78540.2 -> 91999.5
$ws.Cells.Item(46, 13).ClearContents()  # M46: -69769 -> (removed)
$ws.Cells.Item(46, 14).Value = -92461.5  # N46: -79002.2 -> -92461.5
$ws.Cells.Item(113, 8).Value = 459.36365  # H113: 465 -> 459.36365
$ws.Cells.Item(113, 9).Value = 444.25  # I113: 450.14285 -> 444.25
$ws.Cells.Item(113, 11).Value = 1332.75  # K113: 1350.42855 -> 1332.75
$ws.Cells.Item(113, 13).Value = 837.25  # M113: 819.5714499999999 -> 837.25
$ws.Cells.Item(122, 8).Value = 1934.35  # H122: 1974.3684 -> 1934.35
$ws.Cells.Item(122, 9).Value = 2001.5  # I122: 2033.0741 -> 2001.5
$ws.Cells.Item(122, 10).Value = 1777.6666  # J122: 1830.2727 -> 1777.6666
$ws.Cells.Item(122, 11).Value = 6004.5  # K122: 6099.2223 -> 6004.5
$ws.Cells.Item(122, 12).Value = 5332.9998  # L122: 5490.8181 -> 5332.9998
$ws.Cells.Item(122, 13).Value = -3554.5  # M122: -3649.2223 -> -3554.5
$ws.Cells.Item(122, 14).Value = -10232.9998  # N122: -10390.8181 -> -10232.9998
$ws.Cells.Item(134, 8).Value = 91999.5  # H134: 77116.836 -> 91999.5
$ws.Cells.Item(134, 9).Value = 0  # I134: 70000 -> 0
$ws.Cells.Item(134, 10).Value = 91999.5  # J134: 78540.2 -> 91999.5
$ws.Cells.Item(134, 11).Value = 0  # K134: 210000 -> 0
$ws.Cells.Item(134, 12).Value = 275998.5  # L134: 235620.6 -> 275998.5
$ws.Cells.Item(134, 13).ClearContents()  # M134: -207465 -> (removed)
$ws.Cells.Item(134, 14).Value = -281068.5  # N134: -240690.6 -> -281068.5
